# Update the raw-score lookup tables on both age-strata sheets ("060" and "066").
$wb = $excel.ActiveWorkbook

$ws060 = $wb.Worksheets.Item("060")
$ws066 = $wb.Worksheets.Item("066")

# Sheet "060": new COG/EMO raw-score values for ages 1-5
$ws060.Range("B2").Value = 85
$ws060.Range("C2").Value = 90
$ws060.Range("B3").Value = 89
$ws060.Range("C3").Value = 93
$ws060.Range("B4").Value = 97
$ws060.Range("C4").Value = 101
$ws060.Range("B5").Value = 101
$ws060.Range("C5").Value = 110
$ws060.Range("B6").Value = 107
$ws060.Range("C6").Value = 112

# Sheet "066": new COG/EMO raw-score values for ages 1-5
$ws066.Range("B2").Value = 90
$ws066.Range("C2").Value = 87
$ws066.Range("B3").Value = 93
$ws066.Range("C3").Value = 91
$ws066.Range("B4").Value = 101
$ws066.Range("C4").Value = 100
$ws066.Range("B5").Value = 110
$ws066.Range("C5").Value = 106
$ws066.Range("B6").Value = 112
$ws066.Range("C6").Value = 112

# Move the active selection on sheet "066" (the active tab) to C7
$ws066.Activate()
$ws066.Range("C7").Select()
